# Update scripts with new TPM values (Dkk2-Kremen2 LR pair).
# The underlying analysis dropped the "MuSCs" target-cluster rows and
# recomputed the remaining ECs/FAPs -> Dkk2/Kremen2 -> ECs rows, so the
# sheet shrinks from 5 data+header rows (A1:T5) down to 3 (A1:T3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing rows (old rows 4 and 5) - delete bottom-up so the
# row numbers of the rows still to be removed don't shift.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Row 2: ECs -> Dkk2 -> Kremen2 -> ECs (recomputed specificity columns)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dkk2"
$ws.Range("C2").Value = "Kremen2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1717705
$ws.Range("H2").Value = 0.343541
$ws.Range("I2").Value = 0.04063177891664595
$ws.Range("J2").Value = 0.02745976565347561
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.3967905
$ws.Range("N2").Value = 0.793581
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.06815690258024999
$ws.Range("R2").Value = 0.272627610321
$ws.Range("S2").Value = 0.04063177891664595
$ws.Range("T2").Value = 0.02745976565347561

# Row 3: FAPs -> Dkk2 -> Kremen2 -> ECs (recomputed specificity columns)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Dkk2"
$ws.Range("C3").Value = "Kremen2"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.055721
$ws.Range("H3").Value = 12.167163
$ws.Range("I3").Value = 0.9593682210833541
$ws.Range("J3").Value = 0.9725402343465244
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.3967905
$ws.Range("N3").Value = 0.793581
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.6092715634505
$ws.Range("R3").Value = 9.655629380703001
$ws.Range("S3").Value = 0.9593682210833541
$ws.Range("T3").Value = 0.9725402343465244
